$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E) for worker ROXSANNA COROMOTO ESCALONA PEÑA
# (rows 18-24) was listed in descending order (2210 .. 2204). Update it to
# ascending order (2204 .. 2210) to reflect the refreshed EC database.
$periodos = @("2204", "2205", "2206", "2207", "2208", "2209", "2210")
$row = 18
foreach ($p in $periodos) {
    $ws.Range("E$row").Value = $p
    $row++
}
